# Deep sea double count fix
# Updates recalculated landings/status percentages for rows 4, 5, 20, 22, 23
# and refreshes the explanatory footnote in I24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 4 (FAO area 21) ----
$ws.Range("C4").Value = 1.427322557196325
$ws.Range("D4").Value = 15.12284798418123
$ws.Range("E4").Value = 69.0794593258551
$ws.Range("F4").Value = 15.79769268996366
$ws.Range("G4").Value = 84.20230731003633
$ws.Range("H4").Value = 15.79769268996366

# ---- Row 5 (FAO area 27) ----
$ws.Range("B5").Value = 8.08050894
$ws.Range("C5").Value = 7.38183338196252
$ws.Range("D5").Value = 4.053027175415057
$ws.Range("E5").Value = 82.61771444705293
$ws.Range("F5").Value = 13.32925837753201
$ws.Range("G5").Value = 86.67074162246799
$ws.Range("H5").Value = 13.32925837753201

# ---- Row 20 (Deep Sea) ----
# C20 also switches number format to the 3-decimal style (same as used in
# the Sharks row), matching style index 9 in the original workbook.
$ws.Range("C20").NumberFormat = "#,##0.000"
$ws.Range("C20").Value = 0.07802511163319266
$ws.Range("E20").Value = 55.08837186735091
$ws.Range("F20").Value = 44.91162813264909
$ws.Range("G20").Value = 55.08837186735091
$ws.Range("H20").Value = 44.91162813264909

# ---- Row 22 (Sharks) ----
$ws.Range("C22").Value = 0.05364998000000001
$ws.Range("D22").Value = 49.92385831271513
$ws.Range("E22").Value = 38.86139379735091
$ws.Range("F22").Value = 11.21474788993398
$ws.Range("G22").Value = 88.78525211006604
$ws.Range("H22").Value = 11.21474788993398

# ---- Row 23 (Global) ----
$ws.Range("B23").Value = 80.28050343
$ws.Range("C23").Value = 69.45807733661159
$ws.Range("D23").Value = 26.6472356589297
$ws.Range("E23").Value = 48.37613624522654
$ws.Range("F23").Value = 24.97662809584375
$ws.Range("G23").Value = 75.02337190415624
$ws.Range("H23").Value = 24.97662809584375

# ---- Row 24 (Footnote) ----
$noteLines = @(
    "Note: Percent coverage in this sheet does not reflect reported percent coverage. For the reported percent coverage, ",
    "the landings of 'Deep Sea', 'Salmon', and 'Sharks' are incorporated in the FAO major fishing areas ",
    "from which their landings are reported. Thus, percent coverage calculated from this table will slightly different than reported elsewhere. ",
    "Area landings exclude landings from ISSCAAP codes 61, 62, 63, 64, 71, 72, 73, 74, 81, 82, 83, 91, 92, 93, 94, ",
    "except for stocks which have been incorporated in assessment. ",
    "Tuna status/landings have been incorporated into FAO area weighted percentages, so these will appear different ",
    "compared to tables with Tuna category separated."
)
$note = $noteLines -join "`n"

$ws.Range("I24").Value = $note
